# Automatische test-sync: 2025-06-19 19:08:30
# Append a new incoming-mail log row to the "Logs" sheet and refresh the
# "Dashboard" pivot-style summary count for the "Klacht" category.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

$newRow = 39

$logs.Cells.Item($newRow, 1).Value = "Klacht over levering"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Ik ben niet tevreden over mijn bestelling. Ik hoor graag hoe jullie dit oplossen."
$logs.Cells.Item($newRow, 4).Value = "Klacht"
$logs.Cells.Item($newRow, 6).Value = "2025-06-19 19:08:25"
$logs.Cells.Item($newRow, 7).Value = "Nee"

# Grow the conditional-formatting ranges in column D (category) and column G
# (answered y/n) so they keep covering the whole data body, now through the
# freshly appended row.
$logs.Range("D2:D39").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D39"))
$logs.Range("G2:G39").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G39"))

# Update the "Klacht" tally on the Dashboard sheet to reflect the new entry.
$dashboard.Cells.Item(4, 2).Value = 7
